$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F4").Value = "spring core"
$ws.Range("F7").Value = "maven"
$ws.Range("F10").Value = "junit"
$ws.Range("F13").Value = "Spring MVC"

# Match the formatting that Excel carries over from the row's existing
# E-column cell (same fill/font/border) onto the newly added F cells.
$ws.Range("E4").Copy()
$ws.Range("F4").PasteSpecial(-4122)

$ws.Range("E7").Copy()
$ws.Range("F7").PasteSpecial(-4122)

$ws.Range("E10").Copy()
$ws.Range("F10").PasteSpecial(-4122)

$ws.Range("E13").Copy()
$ws.Range("F13").PasteSpecial(-4122)

$excel.CutCopyMode = $false

$ws.Range("F14").Select()
